{"js": "// The document is a date line followed by a 5-column table. The table\n// has 20 rows: 5 \"data\" rows (0-based indices 0, 4, 8, 12, 16) each\n// holding one row's worth of division problems, with 3 blank spacer\n// rows after each data row.\n\n// 1) Update the date line.\nconst dateParagraph = context.document.body.paragraphs.getFirst();\ndateParagraph.insertText(\"2025-08-12 Tuesday\", Word.InsertLocation.replace);\n\n// 2) Update the division problems, addressed by (row, column) so the\n// write-back can't be confused by duplicate equation text elsewhere\n// in the table.\nconst table = context.document.body.tables.getFirst();\n\nconst dataRows = [0, 4, 8, 12, 16];\nconst rowValues = [\n  [\"16\u00f72=\", \"88\u00f73=\", \"45\u00f73=\", \"25\u00f74=\", \"55\u00f73=\"],\n  [\"43\u00f72=\", \"65\u00f73=\", \"75\u00f74=\", \"66\u00f76=\", \"89\u00f76=\"],\n  [\"36\u00f77=\", \"16\u00f72=\", \"10\u00f72=\", \"15\u00f75=\", \"47\u00f76=\"],\n  [\"98\u00f78=\", \"32\u00f75=\", \"52\u00f73=\", \"44\u00f75=\", \"96\u00f76=\"],\n  [\"84\u00f72=\", \"12\u00f79=\", \"71\u00f79=\", \"76\u00f78=\", \"18\u00f77=\"],\n];\n\nfor (let r = 0; r < dataRows.length; r++) {\n  const rowIndex = dataRows[r];\n  const values = rowValues[r];\n  for (let c = 0; c < values.length; c++) {\n    table.getCell(rowIndex, c).value = values[c];\n  }\n}\n\nawait context.sync();\n", "ps1": "# Update the date line at the top of the document.\n$d = $word.ActiveDocument\n$d.Paragraphs.Item(1).Range.Text = \"2025-08-12 Tuesday\"\n\n# Update the division problems laid out in the 5-column table.\n# Each populated table row is followed by three blank spacer rows,\n# so populated rows are Word rows 1, 5, 9, 13, 17 (1-based).\n$t = $d.Tables.Item(1)\n\n$values = @(\n    @(1, 1, \"16\u00f72=\"),\n    @(1, 2, \"88\u00f73=\"),\n    @(1, 3, \"45\u00f73=\"),\n    @(1, 4, \"25\u00f74=\"),\n    @(1, 5, \"55\u00f73=\"),\n    @(5, 1, \"43\u00f72=\"),\n    @(5, 2, \"65\u00f73=\"),\n    @(5, 3, \"75\u00f74=\"),\n    @(5, 4, \"66\u00f76=\"),\n    @(5, 5, \"89\u00f76=\"),\n    @(9, 1, \"36\u00f77=\"),\n    @(9, 2, \"16\u00f72=\"),\n    @(9, 3, \"10\u00f72=\"),\n    @(9, 4, \"15\u00f75=\"),\n    @(9, 5, \"47\u00f76=\"),\n    @(13, 1, \"98\u00f78=\"),\n    @(13, 2, \"32\u00f75=\"),\n    @(13, 3, \"52\u00f73=\"),\n    @(13, 4, \"44\u00f75=\"),\n    @(13, 5, \"96\u00f76=\"),\n    @(17, 1, \"84\u00f72=\"),\n    @(17, 2, \"12\u00f79=\"),\n    @(17, 3, \"71\u00f79=\"),\n    @(17, 4, \"76\u00f78=\"),\n    @(17, 5, \"18\u00f77=\")\n)\n\nforeach ($entry in $values) {\n    $row = $entry[0]\n    $col = $entry[1]\n    $text = $entry[2]\n    $t.Cell($row, $col).Range.Text = $text\n}\n"}
